$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update "Quantity Available" (column K) values for a handful of parts ---
$ws.Range("K2").Value  = 307398
$ws.Range("K5").Value  = 279650
$ws.Range("K8").Value  = 98305
$ws.Range("K23").Value = 1804032
$ws.Range("K25").Value = 209885

# --- Append three new pin-header rows (32-34) to the BOM ---
$newRows = @(
    @{ Row=32; A="SSW-104-01-G-D";     B="Samtec Inc."; C="SAM1208-04-ND";      F="Bulk"; G="Active"; H=1; I="1,16000"; J="1,16 €"; K=2898; L="2 Weeks"; M="CONN RCPT 8POS 0.1 GOLD PCB" },
    @{ Row=33; A="TSW-104-17-G-D";     B="Samtec Inc."; C="SAM1060-04-ND";      F="Bulk"; G="Active"; H=1; I="0,90000"; J="0,90 €"; K=414;  L="2 Weeks"; M="CONN HEADER VERT 8POS 2.54MM" },
    @{ Row=34; A="PH1-10-UA";          B="Adam Tech";   C="2057-PH1-10-UA-ND"; F="Bulk"; G="Active"; H=1; I="0,16000"; J="0,16 €"; K=3470; L="9 Weeks"; M="CONN HEADER VERT 10POS 2.54MM" }
)

foreach ($r in $newRows) {
    $row = $r.Row
    $ws.Cells.Item($row, 1).Value  = $r.A   # Manufacturer Part Number
    $ws.Cells.Item($row, 2).Value  = $r.B   # Manufacturer
    $ws.Cells.Item($row, 3).Value  = $r.C   # Digi-Key Part Number
    $ws.Cells.Item($row, 4).Value  = ""     # Customer Reference
    $ws.Cells.Item($row, 5).Value  = ""     # Reference Designator
    $ws.Cells.Item($row, 6).Value  = $r.F   # Packaging
    $ws.Cells.Item($row, 7).Value  = $r.G   # Part Status
    $ws.Cells.Item($row, 8).Value  = $r.H   # Quantity
    $ws.Cells.Item($row, 9).Value  = $r.I   # Unit Price
    $ws.Cells.Item($row, 10).Value = $r.J   # Extended Price
    $ws.Cells.Item($row, 11).Value = $r.K   # Quantity Available
    $ws.Cells.Item($row, 12).Value = $r.L   # Mfg Std Lead Time
    $ws.Cells.Item($row, 13).Value = $r.M   # Description
    $ws.Cells.Item($row, 14).Value = "ROHS3 Compliant"   # RoHS Status
    $ws.Cells.Item($row, 15).Value = "Lead free"         # Lead Free Status
    $ws.Cells.Item($row, 16).Value = "REACH Unaffected"  # REACH Status
}
